# The commit swaps the raw contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml (the deck's slide-master theme) goes
# from the "Integral" colour scheme to the stock "Office" colour scheme
# (and the reverse for theme2.xml, which backs the notes master).
#
# The font scheme and format scheme blocks are already byte-identical
# between the two theme parts, so the only observable change is the
# clrScheme palette (+ the cosmetic theme/clrScheme "name" attributes).
# We drive that through the PowerPoint colour-scheme object model, which
# is the supported way to rewrite a theme's twelve colour slots.

$p = $ppt.ActivePresentation

function Set-ThemeColors($themeColorScheme, $colors) {
    for ($i = 1; $i -le $colors.Count; $i++) {
        $themeColorScheme.Item($i).RGB = $colors[$i - 1]
    }
}

# RGB() packs as R + G*256 + B*65536 (matches VBA's RGB() macro).
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office" theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in that fixed order.
$officeColors = @(
    (RGB 0x00 0x00 0x00),  # dk1
    (RGB 0xFF 0xFF 0xFF),  # lt1
    (RGB 0x44 0x54 0x6A),  # dk2
    (RGB 0xE7 0xE6 0xE6),  # lt2
    (RGB 0x5B 0x9B 0xD5),  # accent1
    (RGB 0xED 0x7D 0x31),  # accent2
    (RGB 0xA5 0xA5 0xA5),  # accent3
    (RGB 0xFF 0xC0 0x00),  # accent4
    (RGB 0x44 0x72 0xC4),  # accent5
    (RGB 0x70 0xAD 0x47),  # accent6
    (RGB 0x05 0x63 0xC1),  # hlink
    (RGB 0x95 0x4F 0x72)   # folHlink
)

# The slide master's theme is the document theme that PowerPoint's
# object model exposes/edits; apply the swapped-in "Office" palette to
# it there.
$master = $p.SlideMaster
Set-ThemeColors $master.Theme.ThemeColorScheme $officeColors

try { $master.Theme.Name = "Office Theme" } catch {}
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
